$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 117.5
$ws.Range("I33").Value = 120
$ws.Range("J33").Value = 116
$ws.Range("K33").Value = 120
$ws.Range("L33").Value = 116
$ws.Range("M33").Value = 109
$ws.Range("N33").Value = -574

$ws.Range("H51").Value = 5555.4443
$ws.Range("I51").Value = 6399.8
$ws.Range("J51").Value = 4500
$ws.Range("K51").Value = 6399.8
$ws.Range("L51").Value = 4500
$ws.Range("M51").Value = -5915.8
$ws.Range("N51").Value = -5468

$ws.Range("H62").Value = 1941.4166
$ws.Range("I62").Value = 1912.25
$ws.Range("J62").Value = 1999.75
$ws.Range("K62").Value = 1912.25
$ws.Range("L62").Value = 1999.75
$ws.Range("M62").Value = -1288.25
$ws.Range("N62").Value = -3247.75

$ws.Range("H65").Value = 1941.4166
$ws.Range("I65").Value = 1912.25
$ws.Range("J65").Value = 1999.75
$ws.Range("K65").Value = 9561.25
$ws.Range("L65").Value = 9998.75
$ws.Range("M65").Value = -6441.25
$ws.Range("N65").Value = -16238.75

$ws.Range("H111").Value = 11112459
$ws.Range("I111").Value = 14286743
$ws.Range("J111").Value = 2465.5
$ws.Range("K111").Value = 42860229
$ws.Range("L111").Value = 7396.5
$ws.Range("M111").Value = -42857162
$ws.Range("N111").Value = -13530.5

$ws.Range("H116").Value = 11215.154
$ws.Range("I116").Value = 26374.25
$ws.Range("J116").Value = 4477.778
$ws.Range("K116").Value = 26374.25
$ws.Range("L116").Value = 4477.778
$ws.Range("M116").Value = -22932.25
$ws.Range("N116").Value = -11361.778

$ws.Range("H137").Value = 33513.97
$ws.Range("I137").Value = 991.05
$ws.Range("J137").Value = 92646.55
$ws.Range("K137").Value = 2973.15
$ws.Range("L137").Value = 277939.65
$ws.Range("M137").Value = -423.1499999999996
$ws.Range("N137").Value = -283039.65

$ws.Range("H141").Value = 905369.7
$ws.Range("I141").Value = 1038037.2
$ws.Range("K141").Value = 3114111.6
$ws.Range("M141").Value = -3108931.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2906.57
$ws.Range("I32").Value = 2528.602
$ws.Range("J32").Value = 7928.143
$ws.Range("K32").Value = 2528.602
$ws.Range("L32").Value = 7928.143
$ws.Range("M32").Value = -2241.602
$ws.Range("N32").Value = -8502.143

$ws.Range("H88").Value = 3328.3845
$ws.Range("I88").Value = 1551.6666
$ws.Range("J88").Value = 3861.4
$ws.Range("K88").Value = 1551.6666
$ws.Range("L88").Value = 3861.4
$ws.Range("M88").Value = -1145.6666
$ws.Range("N88").Value = -4673.4

$ws.Range("H91").Value = 3328.3845
$ws.Range("I91").Value = 1551.6666
$ws.Range("J91").Value = 3861.4
$ws.Range("K91").Value = 1551.6666
$ws.Range("L91").Value = 3861.4
$ws.Range("M91").Value = -147.6666
$ws.Range("N91").Value = -6669.4

$ws.Range("H122").Value = 1515.6578
$ws.Range("I122").Value = 1439.9688
$ws.Range("K122").Value = 4319.9064
$ws.Range("M122").Value = -1869.9064

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1377.1666
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 1554.3334
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 1554.3334
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -4550.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 150
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -426

$ws.Range("H31").Value = 2653.2144
$ws.Range("I31").Value = 2390.7144
$ws.Range("K31").Value = 2390.7144
$ws.Range("M31").Value = -2095.7144

$ws.Range("H34").Value = 2653.2144
$ws.Range("I34").Value = 2390.7144
$ws.Range("K34").Value = 2390.7144
$ws.Range("M34").Value = -2188.7144

$ws.Range("H99").Value = 1912.0625
$ws.Range("I99").Value = 1589.6
$ws.Range("J99").Value = 2449.5
$ws.Range("K99").Value = 1589.6
$ws.Range("L99").Value = 2449.5
$ws.Range("M99").Value = -91.59999999999991
$ws.Range("N99").Value = -5445.5

$ws.Range("H126").Value = 1912.0625
$ws.Range("I126").Value = 1589.6
$ws.Range("J126").Value = 2449.5
$ws.Range("K126").Value = 4768.799999999999
$ws.Range("L126").Value = 7348.5
$ws.Range("M126").Value = -2298.799999999999
$ws.Range("N126").Value = -12288.5

$ws.Range("H132").Value = 2155.2778
$ws.Range("I132").Value = 1491
$ws.Range("J132").Value = 4480.25
$ws.Range("K132").Value = 4473
$ws.Range("L132").Value = 13440.75
$ws.Range("M132").Value = -1943
$ws.Range("N132").Value = -18500.75

$ws.Range("H134").Value = 2041.1333
$ws.Range("I134").Value = 1907.7097
$ws.Range("J134").Value = 2336.5715
$ws.Range("K134").Value = 5723.1291
$ws.Range("L134").Value = 7009.7145
$ws.Range("M134").Value = -3188.1291
$ws.Range("N134").Value = -12079.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 418.8
$ws.Range("I5").Value = 311.1875
$ws.Range("J5").Value = 849.25
$ws.Range("K5").Value = 933.5625
$ws.Range("L5").Value = 2547.75
$ws.Range("M5").Value = -821.5625
$ws.Range("N5").Value = -2771.75

$ws.Range("H122").Value = 1603.8572
$ws.Range("I122").Value = 735
$ws.Range("J122").Value = 1808.2941
$ws.Range("K122").Value = 6615
$ws.Range("L122").Value = 16274.6469
$ws.Range("M122").Value = -4165
$ws.Range("N122").Value = -21174.6469

$ws.Range("H135").Value = 418.8
$ws.Range("I135").Value = 311.1875
$ws.Range("J135").Value = 849.25
$ws.Range("K135").Value = 2800.6875
$ws.Range("L135").Value = 7643.25
$ws.Range("M135").Value = -265.6875
$ws.Range("N135").Value = -12713.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1114.1578
$ws.Range("I97").Value = 1015.26666
$ws.Range("J97").Value = 1485
$ws.Range("K97").Value = 1015.26666
$ws.Range("L97").Value = 1485
$ws.Range("M97").Value = -519.26666
$ws.Range("N97").Value = -2477

$ws.Range("H126").Value = 1416175.2
$ws.Range("I126").Value = 3088642.2
$ws.Range("J126").Value = 47793.184
$ws.Range("K126").Value = 9265926.600000001
$ws.Range("L126").Value = 143379.552
$ws.Range("M126").Value = -9263456.600000001
$ws.Range("N126").Value = -148319.552

$ws.Range("H132").Value = 653661.8
$ws.Range("I132").Value = 1013448.3
$ws.Range("K132").Value = 3040344.9
$ws.Range("M132").Value = -3037814.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 719.2222
$ws.Range("I100").Value = 496.14285
$ws.Range("K100").Value = 992.2857
$ws.Range("M100").Value = -451.2857
